$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '26.169.30'
Set-TextValue 'E2' '  -6.60%  '
Set-TextValue 'D3' '1.672.27'
Set-TextValue 'E3' '  -4.01%  '
Set-TextValue 'D4' '1.004'
Set-TextValue 'E4' '  +0.19%  '
Set-TextValue 'D5' '218.20'
Set-TextValue 'E5' '  -3.61%  '
Set-TextValue 'D6' '0.5075'
Set-TextValue 'E6' '  -12.21%  '
Set-TextValue 'E7' '  +0.19%  '
Set-TextValue 'D8' '0.2644'
Set-TextValue 'E8' '  -2.41%  '
Set-TextValue 'E9' '  -3.59%  '
Set-TextValue 'D10' '21.53'
Set-TextValue 'E10' '  -6.88%  '
Set-TextValue 'D11' '0.07358'
Set-TextValue 'E11' '  -2.33%  '
Set-TextValue 'D12' '4.558'
Set-TextValue 'E12' '  -3.25%  '
Set-TextValue 'D13' '1.669.50'
Set-TextValue 'E13' '  -4.10%  '
Set-TextValue 'D14' '0.5780'
Set-TextValue 'E14' '  -4.34%  '
Set-TextValue 'D15' '1.896.48'
Set-TextValue 'E15' '  -4.22%  '
Set-TextValue 'D16' '0.000008524'
Set-TextValue 'E16' '  -1.71%  '
Set-TextValue 'D17' '65.06'
Set-TextValue 'E17' '  -12.40%  '
Set-TextValue 'D18' '26.189.70'
Set-TextValue 'E18' '  -6.54%  '
Set-TextValue 'D19' '4.965'
Set-TextValue 'E19' '  -6.68%  '
Set-TextValue 'D20' '1.005'
Set-TextValue 'E20' '  +0.25%  '
Set-TextValue 'D21' '10.80'
Set-TextValue 'E21' '  -4.07%  '
Set-TextValue 'D22' '189.20'
Set-TextValue 'E22' '  -7.73%  '
Set-TextValue 'D23' '6.199'
Set-TextValue 'E23' '  -6.41%  '
Set-TextValue 'E24' '  +0.20%  '
Set-TextValue 'D25' '143.00'
Set-TextValue 'E25' '  -4.55%  '
Set-TextValue 'D26' '7.675'
Set-TextValue 'E26' '  -4.69%  '
Set-TextValue 'D27' '0.1180'
Set-TextValue 'E27' '  -4.11%  '
Set-TextValue 'D28' '15.80'
Set-TextValue 'E28' '  -1.93%  '
Set-TextValue 'D29' '0.05884'
Set-TextValue 'E29' '  -4.73%  '
Set-TextValue 'D30' '1.292'
Set-TextValue 'E30' '  -6.51%  '
Set-TextValue 'D31' '1.320'
Set-TextValue 'E31' '  -5.10%  '
Set-TextValue 'D32' '3.505'
Set-TextValue 'E32' '  -6.25%  '
Set-TextValue 'D33' '3.507'
Set-TextValue 'E33' '  -5.44%  '
Set-TextValue 'D34' '1.660'
Set-TextValue 'E34' '  -0.72%  '
Set-TextValue 'D35' '1.010'
Set-TextValue 'E35' '  -2.40%  '
Set-TextValue 'D36' '0.6008'
Set-TextValue 'E36' '  -5.68%  '
Set-TextValue 'D37' '2.360'
Set-TextValue 'E37' '  -2.48%  '
Set-TextValue 'E38' '  -3.00%  '
Set-TextValue 'E39' '  -3.80%  '
Set-TextValue 'D40' '1.089.85'
Set-TextValue 'E40' '  -3.39%  '
Set-TextValue 'D41' '5.968'
Set-TextValue 'E41' '  -3.70%  '
Set-TextValue 'D42' '0.8592'
Set-TextValue 'E42' '  -1.49%  '
Set-TextValue 'D43' '1.007'
Set-TextValue 'E43' '  +0.30%  '
Set-TextValue 'D44' '99.38'
Set-TextValue 'E44' '  -0.19%  '
Set-TextValue 'D45' '1.821.05'
Set-TextValue 'E45' '  -3.81%  '
Set-TextValue 'E46' '  +1.34%  '
Set-TextValue 'D47' '55.91'
Set-TextValue 'E47' '  -5.76%  '
Set-TextValue 'D48' '1.006'
Set-TextValue 'E48' '  +0.87%  '
Set-TextValue 'D49' '8.061'
Set-TextValue 'E49' '  -2.51%  '
Set-TextValue 'D50' '0.4296'
Set-TextValue 'E50' '  -2.86%  '
Set-TextValue 'D51' '0.05179'
Set-TextValue 'E51' '  -3.71%  '
